$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.719.88'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.02%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.040.30'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.42%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.90%  '

$ws.Range('E6').Value = '  +2.09%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.35'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.384'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.35%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0804'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.82%  '

$ws.Range('E11').Value = '  -0.71%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.343.01'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.69%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.45'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.21%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.06%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.17%  '

$ws.Range('E16').Value = '  +0.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.037.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.46%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.590.23'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.21%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.01%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.63%  '

$ws.Range('E21').Value = '  +1.07%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.77%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.08%  '

$ws.Range('E24').Value = '  +1.11%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.82%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.27'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.20%  '

$ws.Range('E27').Value = '  -0.43%  '

$ws.Range('E28').Value = '  +5.08%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.98%  '

$ws.Range('E30').Value = '  -0.27%  '

$ws.Range('E31').Value = '  +1.19%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.22%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0609'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.87%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.50'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.16%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.01'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.64%  '

$ws.Range('E36').Value = '  -0.30%  '

$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.93'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +11.32%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.24'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.41%  '

$ws.Range('E39').Value = '  +0.18%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.477.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.19%  '

$ws.Range('B41').Value = 'Cronos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0948'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.84%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0216'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.43%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.27%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '95.68'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.53%  '

$ws.Range('E45').Value = '  -0.90%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.15'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +18.03%  '

$ws.Range('E47').Value = '  -2.05%  '

$ws.Range('E48').Value = '  +1.11%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.13%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.52%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.229.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.70%  '
